$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.431.71'
$ws.Range("E2").Value = '  -2.22%  '
$ws.Range("D3").Value = '1.654.94'
$ws.Range("E3").Value = '  -2.01%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.63'
$ws.Range("D5").ClearFormats()
$ws.Range("E6").Value = '  -1.71%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.15'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.261'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -1.30%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0615'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0878'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -0.83%  '
$ws.Range("D12").Value = '1.889.06'
$ws.Range("D13").Value = '1.650.99'
$ws.Range("E13").Value = '  -2.24%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.09'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("E15").Value = '  +2.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.87'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -1.57%  '
$ws.Range("D17").Value = '27.441.78'
$ws.Range("E17").Value = '  -2.02%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.34'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -7.16%  '
$ws.Range("D19").Value = '0.0₃0727'
$ws.Range("E19").Value = '  -2.11%  '
$ws.Range("E20").Value = '  -1.99%  '
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.39'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.31'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -2.72%  '
$ws.Range("E24").Value = '  -1.53%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.72'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.21'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.16%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.94'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.29%  '
$ws.Range("E29").Value = '  -2.28%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0497'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.19'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.36%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.30'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.58%  '
$ws.Range("D33").Value = '1.462.23'
$ws.Range("E33").Value = '  +2.45%  '
$ws.Range("E34").Value = '  -2.66%  '
$ws.Range("E35").Value = '  -4.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.39'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.02%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.912'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.60%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.572'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.33%  '
$ws.Range("E39").Value = '  -1.32%  '
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.46'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.81%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.66'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -5.70%  '
$ws.Range("D45").Value = '1.797.81'
$ws.Range("E45").Value = '  -1.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.782'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -1.83%  '
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.44'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("E49").Value = '  -4.05%  '
$ws.Range("E50").Value = '  -1.56%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.80'
$ws.Range("D51").ClearFormats()
